# TemplateFourCamera.xlsx — make room for a 5th example row.
#
# The sheet currently has 3 header/instruction rows (1-3) followed directly
# by 4 data rows (4-7: pol0/45/90/135 header labels + 3 example file-name
# rows). We insert two blank rows before the data block (pushing the data
# down to rows 6-9) and leave rows 4-5 completely empty, then move the
# selection onto the (now relocated) data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the data block (old row 4 -> new row 6).
$ws.Rows("4:5").Insert()

# Excel's Insert() copies formatting from the row above by default; strip
# that back out so rows 4-5 are genuinely blank (no style, no content),
# matching a plain untouched row.
$ws.Rows("4:5").ClearFormats()
$ws.Rows("4:5").ClearContents()

# Reselect the relocated data block.
$null = $ws.Range("A6:D9").Select()
